# Activities Test data changes - 13 Dec 2023
# Update the Users sheet: replace the attendee name "Drew Koecher" with
# "Ayati Arvind" in cell A2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Range("A2").Value = "Ayati Arvind"
